# "added GDP growth rate" - inserts a new "GDP Growth Rate" data-coverage
# column into the "Data coverage" sheet (sheet1), repurposing the old
# "Quaterly / GDP YoY" column (O) as "GDP Annual / Growth Rate" and adding a
# brand-new "GDP Growth / Rate" column (P), pushing the former P/Q
# (Manufacturing-PMI / OECD Bussiness-Confidence) columns one place right to
# Q/R.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data coverage")

# 1) Insert a new blank column at P - shifts old P->Q and old Q->R, and
#    gives the new column the same per-row formatting as its left neighbour
#    (column O), exactly like Excel's native "Insert Column" command.
$ws.Columns("P").Insert()

# 2) Re-purpose the O column header (row1/row2) text.
$ws.Range("O1").Value = "GDP Annual"
$ws.Range("O2").Value = "Growth Rate"

# 3) Give the freshly inserted P column its header text.
$ws.Range("P1").Value = "GDP Growth"
$ws.Range("P2").Value = "Rate"

# 4) Fill in the new GDP Growth Rate coverage data points. Copy formatting
#    from existing same-styled cells so the fill colours used to flag
#    "year data became available" match the rest of the sheet exactly.

# Countries table style "4" (orange highlight) used for an explicit year.
$ws.Range("J4").Copy()
$ws.Range("P30").PasteSpecial(-4122)
$ws.Range("P30").Value = 2011

$ws.Range("J4").Copy()
$ws.Range("P35").PasteSpecial(-4122)
$ws.Range("P35").Value = 2010

$ws.Range("J4").Copy()
$ws.Range("P53").PasteSpecial(-4122)
$ws.Range("P53").Value = 2005

# Row 57 uses the "to 2021" shared label (style "16", yellow highlight) -
# copy both format and value from N5 which already carries that exact
# combination.
$ws.Range("N5").Copy()
$ws.Range("P57").PasteSpecial(-4122)
$ws.Range("P57").Value = "to 2021"

$excel.CutCopyMode = 0

# 5) Restore view state: the sheet is scrolled down with M32 selected, and
#    the workbook window itself is shifted to the second monitor.
$ws.Range("M32").Select()
$ws.Application.ActiveWindow.ScrollRow = 28
